$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- nChildren count: 4 -> 5 ---
Set-TextValue "B4" "5"

# --- Insert a brand-new data row for the 5th child before the old row 10 ---
# (this shifts the old rows 10,11,12 down to 11,12,13)
$ws.Rows(10).Insert()

# --- Row 6 (child 0): regenerated random data ---
Set-TextValue "B6" "11"
$ws.Range("C6").Value = "Randolph  "
$ws.Range("D6").Value = "Bridgette  "
$ws.Range("E6").Value = "-4,-1"
$ws.Range("F6").Value = "Lenny(father): 0505536740"
Set-TextValue "H6" "14.0"

# --- Row 7 (child 1): regenerated random data ---
Set-TextValue "B7" "13"
$ws.Range("C7").Value = "Fay  "
$ws.Range("D7").Value = "Emilee  "
$ws.Range("E7").Value = "-5,-2"
$ws.Range("F7").Value = "Sheri(mother): 0516797453"
$ws.Range("G7").Value = "7:02:00"
Set-TextValue "H7" "12.0"

# --- Row 8 (child 2): regenerated random data ---
Set-TextValue "B8" "18"
$ws.Range("C8").Value = "Kandis  "
$ws.Range("D8").Value = "Zulma  "
$ws.Range("E8").Value = "-5,-3"
$ws.Range("F8").Value = "Kylie(mother): 0575413269"
$ws.Range("G8").Value = "7:04:00"
Set-TextValue "H8" "10.0"

# --- Row 9 (child 3): regenerated random data ---
Set-TextValue "B9" "12"
$ws.Range("C9").Value = "Frankie  "
$ws.Range("D9").Value = "Flavia  "
$ws.Range("E9").Value = "-3,-4"
$ws.Range("F9").Value = "Cyrus(mother): 0522363358"
$ws.Range("G9").Value = "7:07:00"
Set-TextValue "H9" "7.0"

# --- Row 10 (child 4): brand new row ---
Set-TextValue "A10" "4"
Set-TextValue "B10" "10"
$ws.Range("C10").Value = "Demetra  "
$ws.Range("D10").Value = "Francene  "
$ws.Range("E10").Value = "-2,-3"
$ws.Range("F10").Value = "Dorian(mother): 0534328089"
$ws.Range("G10").Value = "7:09:00"
Set-TextValue "H10" "5.0"

# --- Row 11 (was row 10: school info), start time updated ---
$ws.Range("G11").Value = "7:14:00"

# --- Row 13 (was row 12: time), total time updated ---
Set-TextValue "B13" "14.0"
